$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-02-17 06:39:05"
$wsZhCn.Range("G5").Value = "2016-02-17 06:39:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-02-17 06:39:15"
$wsDeDe.Range("G5").Value = "2016-02-17 06:40:04"
